$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '52.325.11'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.916.31'
$ws.Range('E3').Value = '  +3.69%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '351.64'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.50'
$ws.Range('E6').Value = '  +2.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.560'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.634'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.15'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0866'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.136'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.98'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.80'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.370.46'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  +6.78%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.926.69'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '52.338.74'
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.66'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.33'
$ws.Range('E20').Value = '  +4.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.23'
$ws.Range('E21').Value = '  +3.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0982'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.08'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '271.07'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.80'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.81'
$ws.Range('E26').Value = '  +2.39%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.168'
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.67'
$ws.Range('E29').Value = '  +2.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '37.80'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.27'
$ws.Range('E31').Value = '  +10.65%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.50'
$ws.Range('E33').Value = '  +5.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0964'
$ws.Range('E34').Value = '  +10.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '53.18'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.33'
$ws.Range('E38').Value = '  +5.83%  '
$ws.Range('E39').Value = '  +16.43%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.83'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '23.84'
$ws.Range('E42').Value = '  +8.13%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.62'
$ws.Range('E45').Value = '  +6.47%  '
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.57'
$ws.Range('E47').Value = '  +5.39%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.197.24'
$ws.Range('E48').Value = '  +4.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.265'
$ws.Range('E49').Value = '  +24.35%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0337'
$ws.Range('E50').Value = '  +10.65%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.965'
$ws.Range('E51').Value = '  +3.88%  '
